# Apply updated dSF ("F" column) values per the "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -7
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = -5
$ws.Range("F15").Value = -2
